$wb = $excel.ActiveWorkbook

# --- Update the summary text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.35 = 30221.32 pesos`n✅ 30221.32 pesos = 7.34 = 956.23 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 136
$wsTasas.Range("O10").Value = 4110.1
$wsTasas.Range("N12").Value = 4120
$wsTasas.Range("O12").Value = 130.36
